$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like
# "27.523.67" or "322.77" are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.523.67"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.744.43"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "322.77"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.4452"
$ws.Range("E7").Value = "  +4.79%  "
$ws.Range("D8").Value = "0.3522"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").Value = "0.07410"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "41.56"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "5.904"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "7.090"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "1.743.01"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "91.56"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "0.06381"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "16.88"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "5.726"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("D23").Value = "27.553.82"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "2.096"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "160.72"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "1.942.89"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "125.27"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "2.033"
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("D32").Value = "0.09081"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("D33").Value = "3.653"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "5.368"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("D35").Value = "0.02275"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("D37").Value = "0.06033"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "0.2064"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").Value = "4.895"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "0.6231"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "1.375"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "7.715"
$ws.Range("D44").Value = "13.18"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "3.699"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "0.5797"
$ws.Range("D47").Value = "122.14"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "0.06841"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "1.116"
$ws.Range("E50").Value = "  -4.19%  "
$ws.Range("D51").Value = "71.47"
$ws.Range("E51").Value = "  -2.18%  "

# Restore the original (default) cell style so no style/number-format
# metadata is left behind on the price cells.
$priceRange.Style = "Normal"
